# Meeting20150425/Social Media.pptx - update Followers table:
#   - new table style
#   - Twitter followers 144 -> 149
#   - Google+ followers 40 -> 41

$p = $ppt.ActivePresentation

# Locate the "Followers" slide: it is the slide that holds a table shape.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table

            # Apply the new table style.
            $tbl.ApplyStyle("{BA884226-E885-4FA7-BA9F-EDDA90F8EC01}")

            # Walk every cell and update the follower counts in place.
            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cell = $tbl.Cell($r, $c)
                    $txt = $cell.Shape.TextFrame.TextRange.Text
                    if ($txt -eq "144") {
                        $cell.Shape.TextFrame.TextRange.Text = "149"
                    } elseif ($txt -eq "40") {
                        $cell.Shape.TextFrame.TextRange.Text = "41"
                    }
                }
            }
        }
    }
}
